$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Target cluster (column D) and Sending cluster (column A) labels ---
$ws.Range("A5").Value2 = "FAPs"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("D7").Value2 = "MuSCs"

# --- Update recomputed TPM-normalized metric values ---
# Row 2
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 0.7103813333333333
$ws.Range("H2").Value2 = 2.131144
$ws.Range("I2").Value2 = 0.7576743564291667
$ws.Range("J2").Value2 = 0.7576743564291667
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 1.414593
$ws.Range("N2").Value2 = 4.243779
$ws.Range("O2").Value2 = 0.3478871232761722
$ws.Range("P2").Value2 = 0.3478871232761722
$ws.Range("Q2").Value2 = 1.004900461464
$ws.Range("R2").Value2 = 9.044104153175999
$ws.Range("S2").Value2 = 0.2635851522382679
$ws.Range("T2").Value2 = 0.263585152238268

# Row 3
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 0.7103813333333333
$ws.Range("H3").Value2 = 2.131144
$ws.Range("I3").Value2 = 0.7576743564291667
$ws.Range("J3").Value2 = 0.7576743564291667
$ws.Range("O3").Value2 = 0.5748520910875596
$ws.Range("P3").Value2 = 0.5748520910875596
$ws.Range("Q3").Value2 = 1.660507368503111
$ws.Range("R3").Value2 = 14.944566316528
$ws.Range("S3").Value2 = 0.4355506881567275
$ws.Range("T3").Value2 = 0.4355506881567275

# Row 4
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 0.7103813333333333
$ws.Range("H4").Value2 = 2.131144
$ws.Range("I4").Value2 = 0.7576743564291667
$ws.Range("J4").Value2 = 0.7576743564291667
$ws.Range("M4").Value2 = 0.314161
$ws.Range("N4").Value2 = 0.942483
$ws.Range("O4").Value2 = 0.07726078563626818
$ws.Range("P4").Value2 = 0.07726078563626819
$ws.Range("Q4").Value2 = 0.2231741100613333
$ws.Range("R4").Value2 = 2.008566990552
$ws.Range("S4").Value2 = 0.0585385160341713
$ws.Range("T4").Value2 = 0.0585385160341713

# Row 5
$ws.Range("E5").Value2 = 2
$ws.Range("F5").Value2 = 0.6666666666666666
$ws.Range("G5").Value2 = 0.2272
$ws.Range("H5").Value2 = 0.6816
$ws.Range("I5").Value2 = 0.2423256435708333
$ws.Range("J5").Value2 = 0.2423256435708333
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 1.414593
$ws.Range("N5").Value2 = 4.243779
$ws.Range("O5").Value2 = 0.3478871232761722
$ws.Range("P5").Value2 = 0.3478871232761722
$ws.Range("Q5").Value2 = 0.3213955296
$ws.Range("R5").Value2 = 2.8925597664
$ws.Range("S5").Value2 = 0.08430197103790425
$ws.Range("T5").Value2 = 0.08430197103790427

# Row 6
$ws.Range("I6").Value2 = 0.2423256435708333
$ws.Range("J6").Value2 = 0.2423256435708333
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 2.337487333333333
$ws.Range("N6").Value2 = 7.012461999999999
$ws.Range("O6").Value2 = 0.5748520910875596
$ws.Range("P6").Value2 = 0.5748520910875596
$ws.Range("Q6").Value2 = 0.5310771221333332
$ws.Range("R6").Value2 = 4.779694099199999
$ws.Range("S6").Value2 = 0.1393014029308322
$ws.Range("T6").Value2 = 0.1393014029308322

# Row 7
$ws.Range("I7").Value2 = 0.2423256435708333
$ws.Range("J7").Value2 = 0.2423256435708333
$ws.Range("M7").Value2 = 0.314161
$ws.Range("N7").Value2 = 0.942483
$ws.Range("O7").Value2 = 0.07726078563626818
$ws.Range("P7").Value2 = 0.07726078563626819
$ws.Range("Q7").Value2 = 0.07137737919999999
$ws.Range("R7").Value2 = 0.6423964127999999
$ws.Range("S7").Value2 = 0.01872226960209688
$ws.Range("T7").Value2 = 0.01872226960209689

# --- Remove now-empty rows 8 and 9 (Resolving-Mac target cluster dropped) ---
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()

